$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new aggregated text for row 2 (Figure of Destiny)
$figureText = "('Figure of Destiny', ['{R/W}', 'Creature " + [char]0x2014 + " Kithkin', '{R/W}: Figure of Destiny becomes a Kithkin Spirit with base power and toughness 2/2.', '{R/W}{R/W}{R/W}: If Figure of Destiny is a Spirit, it becomes a Kithkin Spirit Warrior with base power and toughness 4/4.', '{R/W}{R/W}{R/W}{R/W}{R/W}{R/W}: If Figure of Destiny is a Warrior, it becomes a Kithkin Spirit Warrior Avatar with base power and toughness 8/8, flying, and first strike.', '1/1'])"

# Build the new aggregated text for row 3 (Overbeing of Myth)
$overbeingText = "('Overbeing of Myth', ['{G/U}{G/U}{G/U}{G/U}{G/U}', 'Creature " + [char]0x2014 + " Spirit Avatar', 'Overbeing of Myth" + [char]0x2019 + "s power and toughness are each equal to the number of cards in your hand.', 'At the beginning of your draw step, draw an additional card.', '*/*'])"

# Update A2 and A3 with the combined values
$ws.Range("A2").Value = $figureText
$ws.Range("A3").Value = $overbeingText

# Delete rows 4 through 14, which are no longer needed
$ws.Range("A4:A14").EntireRow.Delete()
